$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Clean up the "Fuel Type" labels in column A: the values were stored with
#    a leading two-space indent ("  Black coal" etc). Re-write them without
#    the indent. This also causes the now-unused indented shared strings to
#    drop out of xl/sharedStrings.xml, with the clean names appended in row
#    order - matching how Excel recompacts the shared-string table.
# ---------------------------------------------------------------------------
$ws.Range("A2").Value  = "Black coal"
$ws.Range("A3").Value  = "Brown coal"
$ws.Range("A4").Value  = "Natural gas"
$ws.Range("A5").Value  = "Oil products"
$ws.Range("A6").Value  = "Other non-renewable"
$ws.Range("A7").Value  = "Bagasse, wood"
$ws.Range("A8").Value  = "Biogas"
$ws.Range("A9").Value  = "Wind"
$ws.Range("A10").Value = "Hydro"
$ws.Range("A11").Value = "Large-scale solar PV"
$ws.Range("A12").Value = "Small-scale solar PV"
$ws.Range("A13").Value = "Geothermal"

# ---------------------------------------------------------------------------
# 2. Add the new 2022 data column (AJ), one column past the existing 1989-2021
#    series (columns C .. AI). The header is a literal year value (not a
#    "+1" formula like the rest of the header row).
# ---------------------------------------------------------------------------
$ws.Range("AJ1").Value = 2022

$ws.Range("AJ2").Value  = 96173.858193717693
$ws.Range("AJ3").Value  = 31459.125523381801
$ws.Range("AJ4").Value  = 48865.237380704799
$ws.Range("AJ5").Value  = 4864.3971483286441
$ws.Range("AJ6").Value  = 0
$ws.Range("AJ7").Value  = 1687.7196683510799
$ws.Range("AJ8").Value  = 1405.0737606722
$ws.Range("AJ9").Value  = 31384.873888888898
$ws.Range("AJ10").Value = 16666.054444444399
$ws.Range("AJ11").Value = 16059.8195145678
$ws.Range("AJ12").Value = 25908.682159
$ws.Range("AJ13").Value = 0

# Give the new figures the same look as the rest of the pasted-in series:
# one decimal place, right aligned, Arial/black (AJ2:AJ5, AJ7:AJ13), built
# up on a scratch cell and copied over with Paste Special so every cell
# lands on a single shared style instead of one-off per-property styles.
$scratch1 = $ws.Range("ZZ1")
$scratch1.NumberFormat = "#,##0.0"
$scratch1.Font.Color = 0
$scratch1.Font.Name = "Arial"
$scratch1.HorizontalAlignment = -4152
$scratch1.Copy()
$ws.Range("AJ2:AJ5").PasteSpecial(-4122)
$ws.Range("AJ7:AJ13").PasteSpecial(-4122)
$scratch1.Clear()

# Row 6 ("Other non-renewable") came from a separately-sourced paste: same
# number format and font, but unlocked and without the explicit color/
# alignment override.
$scratch2 = $ws.Range("ZZ2")
$scratch2.NumberFormat = "#,##0.0"
$scratch2.Font.Name = "Arial"
$scratch2.Locked = $false
$scratch2.Copy()
$ws.Range("AJ6").PasteSpecial(-4122)
$scratch2.Clear()

$excel.CutCopyMode = $false

# Give the new column the same display width as the rest of the year columns.
$ws.Columns.Item(36).ColumnWidth = $ws.Columns.Item(6).ColumnWidth

# ---------------------------------------------------------------------------
# 3. View state: the sheet had scrolled so column K was pinned at the left
#    edge with AD25 selected; restore the default scroll position and select
#    F18 instead.
# ---------------------------------------------------------------------------
$ws.Range("F18").Select()

# Restore the window size/position recorded for the workbook.
$win = $excel.ActiveWindow
$win.Left = 19090
$win.Top = -1330
$win.Width = 25820
$win.Height = 13900
